$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 06:32"

# Peru row (row 10) - refreshed case counts
$ws.Range("B10").Value = 483133
$ws.Range("D10").Value = 329404
$ws.Range("E10").Value = 132453
$ws.Range("H10").Value = 21276

# Rows 185/186: Bermudas and Belice swap places (label + data move together)
$ws.Range("A185").Value = "Belice"
$ws.Range("B185").Value = 177
$ws.Range("C185").Value = 23
$ws.Range("D185").Value = 32
$ws.Range("E185").Value = 143
$ws.Range("H185").Value = 2

$ws.Range("A186").Value = "Bermudas"
$ws.Range("B186").Value = 158
$ws.Range("D186").Value = 144
$ws.Range("E186").Value = 5
$ws.Range("H186").Value = 9

# Rows 202/203: Timor Oriental and Santa Lucia swap places (data identical, only labels swap)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# Rows 213/214: Montserrat and Islas Malvinas swap places (label + data move together)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
